# Update cryptocurrency price/volume data per Wed Nov 15 03:51:28 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.436.42"
$ws.Range("E2").Value = "'  -2.90%  "
$ws.Range("D3").Value = "'1.972.32"
$ws.Range("E3").Value = "'  -4.10%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'244.56"
$ws.Range("E5").Value = "'  +0.95%  "
$ws.Range("E6").Value = "'  -4.40%  "
$ws.Range("B7").Value = "'Solana"
$ws.Range("C7").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'56.78"
$ws.Range("E7").Value = "'  +4.20%  "
$ws.Range("B8").Value = "'USDC"
$ws.Range("C8").Value = "'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("D9").Value = "'58.15"
$ws.Range("E9").Value = "'  -0.75%  "
$ws.Range("E10").Value = "'  -0.92%  "
$ws.Range("D11").Value = "'0.0731"
$ws.Range("E11").Value = "'  -2.56%  "
$ws.Range("E12").Value = "'  -2.86%  "
$ws.Range("D13").Value = "'0.940"
$ws.Range("E13").Value = "'  +3.94%  "
$ws.Range("D14").Value = "'14.25"
$ws.Range("E14").Value = "'  -3.05%  "
$ws.Range("D15").Value = "'2.262.27"
$ws.Range("E15").Value = "'  -4.06%  "
$ws.Range("D16").Value = "'5.25"
$ws.Range("E16").Value = "'  -2.46%  "
$ws.Range("D17").Value = "'1.968.68"
$ws.Range("E17").Value = "'  -4.29%  "
$ws.Range("D18").Value = "'17.45"
$ws.Range("E18").Value = "'  +4.47%  "
$ws.Range("D19").Value = "'35.382.62"
$ws.Range("E19").Value = "'  -2.92%  "
$ws.Range("D20").Value = "'71.46"
$ws.Range("E20").Value = "'  -0.84%  "
$ws.Range("D21").Value = "'0.0₃0839"
$ws.Range("E21").Value = "'  -2.14%  "
$ws.Range("D22").Value = "'232.07"
$ws.Range("E22").Value = "'  -2.35%  "
$ws.Range("D23").Value = "'5.11"
$ws.Range("E23").Value = "'  -2.20%  "
$ws.Range("E24").Value = "'  +0.05%  "
$ws.Range("D25").Value = "'2.54"
$ws.Range("E25").Value = "'  +19.46%  "
$ws.Range("E26").Value = "'  -2.49%  "
$ws.Range("D27").Value = "'163.94"
$ws.Range("E27").Value = "'  +0.01%  "
$ws.Range("E28").Value = "'  -3.15%  "
$ws.Range("D29").Value = "'19.04"
$ws.Range("E29").Value = "'  -5.16%  "
$ws.Range("E30").Value = "'  -2.82%  "
$ws.Range("D31").Value = "'4.84"
$ws.Range("E31").Value = "'  -3.88%  "
$ws.Range("E32").Value = "'  -6.77%  "
$ws.Range("B33").Value = "'Kaspa"
$ws.Range("C33").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.0918"
$ws.Range("E33").Value = "'  +11.77%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0588"
$ws.Range("E34").Value = "'  -0.78%  "
$ws.Range("E35").Value = "'  -3.60%  "
$ws.Range("E36").Value = "'  +7.43%  "
$ws.Range("E37").Value = "'  -0.06%  "
$ws.Range("D38").Value = "'1.76"
$ws.Range("E38").Value = "'  -3.86%  "
$ws.Range("E39").Value = "'  +5.18%  "
$ws.Range("E40").Value = "'  -2.18%  "
$ws.Range("E41").Value = "'  +1.27%  "
$ws.Range("E42").Value = "'  -2.56%  "
$ws.Range("E43").Value = "'  -2.12%  "
$ws.Range("D44").Value = "'90.76"
$ws.Range("E44").Value = "'  -3.10%  "
$ws.Range("B45").Value = "'Cronos"
$ws.Range("C45").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0881"
$ws.Range("E45").Value = "'  -2.58%  "
$ws.Range("B46").Value = "'Maker"
$ws.Range("C46").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'1.368.93"
$ws.Range("E46").Value = "'  -1.56%  "
$ws.Range("D47").Value = "'15.74"
$ws.Range("E47").Value = "'  -0.20%  "
$ws.Range("E48").Value = "'  -0.84%  "
$ws.Range("E49").Value = "'  +0.85%  "
$ws.Range("D50").Value = "'45.96"
$ws.Range("E50").Value = "'  +1.65%  "
$ws.Range("D51").Value = "'3.62"
$ws.Range("E51").Value = "'  +9.56%  "
